$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Output $ws.GetType()
try { Write-Output $ws.ListObjects.Count } catch { Write-Output "no ListObjects: $_" }
